$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, copying the formatting used by the
# other header cells (e.g. G1) so it gets the same bold/centered/bordered style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the value for the new column in row 2 (plain numeric cell, no special style)
$ws.Range("H2").Value = 0
